$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 3: Krankenversicherung (no longer "Arbeitgeberbeitrag ...")
$ws.Range("A3").Value = "Krankenversicherung in Prozent"
$ws.Range("B3").Value = 13

# Row 4: Arbeitgeberbeitrag Rentenversicherung in Prozent
$ws.Range("B4").Value = 15

# Row 5: Arbeitnehmerbeitrag Rentenversicherung in Prozent
$ws.Range("B5").Value = 3.6

# Row 7: U2-Umlage in Prozent
$ws.Range("B7").Value = 0.24

# Row 10: Eintragsdatum value update
$ws.Range("B10").Value = "01.01.2024"
